# Book1.xlsx edit: "changed boo1 .added 1 row"
# Appends a new row (d / 500) to Sheet1, right after the existing a/b/c rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "d"
$ws.Range("B5").Value = 500

# Mirror the author's saved selection state (active cell on the new row).
$ws.Range("B5").Select()
